$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Move the existing ExpPoints column (C) to G, then write the new header
# columns (WIN, TOP2, TOP4, RELEGATION) into C:F.

# Shift column C (ExpPoints) to column G by inserting 4 blank columns
# before it.
$ws.Range("C1:F1").EntireColumn.Insert()

# New header labels for the inserted columns.
$ws.Range("C1").Value = "WIN"
$ws.Range("D1").Value = "TOP2"
$ws.Range("E1").Value = "TOP4"
$ws.Range("F1").Value = "RELEGATION"

# The inserted columns are left as empty text cells for rows 2-19
# (placeholders for future Monte Carlo simulation percentages). A lone
# "'" forces an empty-text cell (Excel's text-prefix quote) instead of
# a truly blank one; re-applying the Normal style afterwards drops the
# quote-prefix formatting flag again so the cell stays unstyled.
$fillRange = $ws.Range("C2:F19")
$fillRange.Value = "'"
$fillRange.Style = "Normal"
